# Revert "Deploying to gh-pages from @ ncpi-fhir/ncpi-fhir-ig-v0.2@4358b16"
# Reverts URL/version/date/system-URI metadata strings back to prior values.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig/ValueSet/condition-inheritance-vs"
$wsMeta.Range("B3").Value = "0.1.0"
$wsMeta.Range("B8").Value = "2021-12-13T19:24:22+00:00"

$wsInclude = $wb.Worksheets.Item("Include from Condition Inheri")
$wsInclude.Range("B4").Value = "https://ncpi-fhir.github.io/ncpi-fhir-ig/CodeSystem/ConditionInheritanceMode"
